# Fix logistic regression result:
# Insert "Logistic Regression" as the first model row (row 2), pushing the
# other models down by one row, and add a new "Configuration" column (B)
# populated with "NULL" for every model row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A - Model names (Logistic Regression now leads the list)
$ws.Range("A2").Value = "Logistic Regression"
$ws.Range("A3").Value = "Multinomial Naive Bayes"
$ws.Range("A4").Value = "Support Vector Machines"
$ws.Range("A5").Value = "Decision Tree"
$ws.Range("A6").Value = "Random Forest"

# Column B - new Configuration values
$ws.Range("B2").Value = "NULL"
$ws.Range("B3").Value = "NULL"
$ws.Range("B4").Value = "NULL"
$ws.Range("B5").Value = "NULL"
$ws.Range("B6").Value = "NULL"

# Update the active cell / selection to C2
$ws.Range("C2").Select()
